$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet holds a weekly time series of "Cebollín" (green onion) price
# records in rows 2-126 (row 1 is the header). A new weekly record is being
# added. It is inserted above the current row 106, which pushes the
# existing rows 106-126 down to 107-127 (matching the diff, where every
# row's data from 106 on is replaced by the row above's former content,
# and the table grows from A1:R126 to A1:R127).
$ws.Rows.Item(106).Insert()

# After the insert, row 107 holds what used to be row 106's data. Copy it
# into the freshly inserted (blank) row 106 so every "template" column
# (Mercado, Región, Codreg, Categoría ID, Categoría, Variedad, Calidad,
# Precio mínimo/máximo, Unidad de comercialización, Origen, Kg o Unidades,
# Clasificación, etc.) is populated exactly like the surrounding rows.
$ws.Rows.Item(107).Copy()
$ws.Rows.Item(106).PasteSpecial()

# Finally, fill in the two values that are actually new for this week's
# record: the date (Fecha) and the volume (Volumen).
$ws.Range("D106").Value = 44505
$ws.Range("J106").Value = 3100
